$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above current row 2 (shifts old row2->row3, row3->row4)
$ws.Rows.Item(2).Insert()

# New row 2: Creazione_2 / REGIONE_CAMPANIA event
$ws.Cells.Item(2, 1).Value = "Creazione_2"
$ws.Cells.Item(2, 2).Value = "REGIONE_CAMPANIA"
$ws.Cells.Item(2, 3).Value = "PRZMGV95D49F839P^^^&2.16.840.1.113883.2.9.4.3.2&ISO"
$ws.Cells.Item(2, 4).Value = "2.16.840.1.113883.2.9.2.120.4.4.b0f3ffcf25ce2aafc7dc901e2febc51f43837f4ca0fe3b6d1b02194e9047b6db.52d02c742a^^^^urn:ihe:iti:xdw:2013:workflowInstanceId"
$ws.Cells.Item(2, 5).Value = "2.16.840.1.113883.2.9.2.110.4.4^UAT_GTW_ID1721656896931"
$ws.Cells.Item(2, 6).Value = "22-07-2024:16:01:38"

# Row 3 (previously row 2, Creazione_1 / REGIONE_CAMPANIA): refresh workflow id, gateway id, timestamp
$ws.Cells.Item(3, 1).Value = "Creazione_1"
$ws.Cells.Item(3, 2).Value = "REGIONE_CAMPANIA"
$ws.Cells.Item(3, 3).Value = "NGNVCN92S19L259C^^^&2.16.840.1.113883.2.9.4.3.2&ISO"
$ws.Cells.Item(3, 4).Value = "2.16.840.1.113883.2.9.2.120.4.4.b0f3ffcf25ce2aafc7dc901e2febc51f43837f4ca0fe3b6d1b02194e9047b6db.4faa8a94e2^^^^urn:ihe:iti:xdw:2013:workflowInstanceId"
$ws.Cells.Item(3, 5).Value = "2.16.840.1.113883.2.9.2.110.4.4^UAT_GTW_ID1721656887136"
$ws.Cells.Item(3, 6).Value = "22-07-2024:16:01:28"

# Row 4 (previously row 3, Creazione_0 / REGIONE_LAZIO): refresh workflow id, gateway id, timestamp
$ws.Cells.Item(4, 1).Value = "Creazione_0"
$ws.Cells.Item(4, 2).Value = "REGIONE_LAZIO"
$ws.Cells.Item(4, 3).Value = "MRCLSN97C14H501J^^^&2.16.840.1.113883.2.9.4.3.2&ISO"
$ws.Cells.Item(4, 4).Value = "2.16.840.1.113883.2.9.2.120.4.4.b0f3ffcf25ce2aafc7dc901e2febc51f43837f4ca0fe3b6d1b02194e9047b6db.592d0019d3^^^^urn:ihe:iti:xdw:2013:workflowInstanceId"
$ws.Cells.Item(4, 5).Value = "2.16.840.1.113883.2.9.2.110.4.4^UAT_GTW_ID1721656876332"
$ws.Cells.Item(4, 6).Value = "22-07-2024:16:01:17"
